# Updates the five-row division-fact table: each of the 25 populated
# cells (rows 1, 5, 9, 13, 17 of the single table; 5 cells each) gets its
# "NN÷N=" text replaced with a new expression. Cells are addressed by
# (row, column) rather than by text search, because a couple of the new
# values ("92÷9=") coincide with other cells' pre-edit text, which would
# make a naive sequential Find/Replace ambiguous.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Col = 1; Old = "40÷8="; New = "49÷3=" },
    @{ Row = 1;  Col = 2; Old = "97÷9="; New = "87÷5=" },
    @{ Row = 1;  Col = 3; Old = "14÷4="; New = "10÷8=" },
    @{ Row = 1;  Col = 4; Old = "27÷5="; New = "11÷6=" },
    @{ Row = 1;  Col = 5; Old = "62÷2="; New = "15÷2=" },

    @{ Row = 5;  Col = 1; Old = "26÷5="; New = "10÷5=" },
    @{ Row = 5;  Col = 2; Old = "96÷2="; New = "33÷5=" },
    @{ Row = 5;  Col = 3; Old = "85÷9="; New = "79÷8=" },
    @{ Row = 5;  Col = 4; Old = "50÷9="; New = "87÷9=" },
    @{ Row = 5;  Col = 5; Old = "26÷4="; New = "86÷6=" },

    @{ Row = 9;  Col = 1; Old = "41÷7="; New = "52÷4=" },
    @{ Row = 9;  Col = 2; Old = "11÷8="; New = "75÷5=" },
    @{ Row = 9;  Col = 3; Old = "71÷2="; New = "15÷3=" },
    @{ Row = 9;  Col = 4; Old = "36÷2="; New = "57÷7=" },
    @{ Row = 9;  Col = 5; Old = "69÷5="; New = "19÷9=" },

    @{ Row = 13; Col = 1; Old = "90÷8="; New = "14÷5=" },
    @{ Row = 13; Col = 2; Old = "34÷6="; New = "93÷5=" },
    @{ Row = 13; Col = 3; Old = "92÷9="; New = "33÷7=" },
    @{ Row = 13; Col = 4; Old = "58÷9="; New = "67÷8=" },
    @{ Row = 13; Col = 5; Old = "39÷2="; New = "92÷9=" },

    @{ Row = 17; Col = 1; Old = "66÷6="; New = "11÷4=" },
    @{ Row = 17; Col = 2; Old = "20÷3="; New = "92÷9=" },
    @{ Row = 17; Col = 3; Old = "24÷3="; New = "10÷7=" },
    @{ Row = 17; Col = 4; Old = "96÷4="; New = "94÷8=" },
    @{ Row = 17; Col = 5; Old = "49÷8="; New = "21÷7=" }
)

foreach ($ch in $changes) {
    $cell = $t.Cell($ch.Row, $ch.Col)
    $current = $cell.Range.Text
    if ($current -notmatch [regex]::Escape($ch.Old)) {
        throw "Cell ($($ch.Row),$($ch.Col)) expected '$($ch.Old)' but found '$current'"
    }
    $cell.Range.Text = $ch.New
}
